# Auto-generated edit script: update cryptos list values (prices & 1h volume %) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.747.14'
$ws.Range('E2').Value = '  -5.01%  '
$ws.Range('D3').Value = '3.217.07'
$ws.Range('E3').Value = '  -8.28%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '596.62'
$ws.Range('E5').Value = '  -1.33%  '
$ws.Range('D6').Value = '152.18'
$ws.Range('E6').Value = '  -12.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.207.77'
$ws.Range('E8').Value = '  -8.24%  '
$ws.Range('E9').Value = '  -10.38%  '
$ws.Range('E10').Value = '  -10.84%  '
$ws.Range('D11').Value = '6.63'
$ws.Range('E11').Value = '  -7.79%  '
$ws.Range('D12').Value = '0.498'
$ws.Range('E12').Value = '  -15.00%  '
$ws.Range('D13').Value = '39.23'
$ws.Range('E13').Value = '  -14.88%  '
$ws.Range('E14').Value = '  -11.24%  '
$ws.Range('D15').Value = '3.739.18'
$ws.Range('E15').Value = '  -8.26%  '
$ws.Range('D16').Value = '66.804.27'
$ws.Range('E16').Value = '  -4.84%  '
$ws.Range('D17').Value = '3.218.83'
$ws.Range('E17').Value = '  -8.35%  '
$ws.Range('E18').Value = '  -4.36%  '
$ws.Range('D19').Value = '534.39'
$ws.Range('E19').Value = '  -13.04%  '
$ws.Range('D20').Value = '7.19'
$ws.Range('E20').Value = '  -13.12%  '
$ws.Range('D21').Value = '15.02'
$ws.Range('E21').Value = '  -14.22%  '
$ws.Range('D22').Value = '0.765'
$ws.Range('E22').Value = '  -12.97%  '
$ws.Range('D23').Value = '7.95'
$ws.Range('E23').Value = '  -12.38%  '
$ws.Range('E24').Value = '  -10.58%  '
$ws.Range('D25').Value = '86.31'
$ws.Range('E25').Value = '  -13.03%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').Value = '3.22'
$ws.Range('E27').Value = '  -13.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.20'
$ws.Range('E28').Value = '  -14.26%  '
$ws.Range('D29').Value = '8.18'
$ws.Range('E29').Value = '  -9.29%  '
$ws.Range('D30').Value = '29.53'
$ws.Range('E30').Value = '  -13.27%  '
$ws.Range('E31').Value = '  -9.59%  '
$ws.Range('E32').Value = '  -8.88%  '
$ws.Range('D33').Value = '547.79'
$ws.Range('E33').Value = '  -15.88%  '
$ws.Range('D35').Value = '5.74'
$ws.Range('E35').Value = '  -15.68%  '
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.40'
$ws.Range('E37').Value = '  -5.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0430'
$ws.Range('E38').Value = '  -9.65%  '
$ws.Range('E39').Value = '  -13.21%  '
$ws.Range('D40').Value = '9.35'
$ws.Range('E40').Value = '  -12.75%  '
$ws.Range('D41').Value = '0.125'
$ws.Range('E41').Value = '  -12.07%  '
$ws.Range('D42').Value = '2.71'
$ws.Range('E42').Value = '  -22.72%  '
$ws.Range('D43').Value = '2.924.52'
$ws.Range('E43').Value = '  -12.83%  '
$ws.Range('D44').Value = '0.265'
$ws.Range('E44').Value = '  -14.29%  '
$ws.Range('E45').Value = '  -20.33%  '
$ws.Range('E46').Value = '  -15.91%  '
$ws.Range('D47').Value = '26.62'
$ws.Range('E47').Value = '  -16.66%  '
$ws.Range('E48').Value = '  -15.07%  '
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').Value = '0.114'
$ws.Range('E50').Value = '  -11.73%  '
$ws.Range('D51').Value = '120.11'
$ws.Range('E51').Value = '  -9.48%  '
